# Auto-generated edit script applying the commit diff to cryptos.xlsx
# Updates symbol list (coin rotation rows 7-18), price/volume figures, and
# the "Hora" (hour) column from 10 -> 11 across all data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether it must be forced to
# text (numeric-looking strings like prices/percentages/hours need a
# leading quote so Excel stores them as text, matching the source data).
$edits = @(
    @{Cell='D2'; Value='261.11'; AsText=$true}
    @{Cell='E2'; Value='1.84%'; AsText=$true}
    @{Cell='G2'; Value='11'; AsText=$true}
    @{Cell='D3'; Value='27.34'; AsText=$true}
    @{Cell='E3'; Value='1.89%'; AsText=$true}
    @{Cell='G3'; Value='11'; AsText=$true}
    @{Cell='D4'; Value='4.710'; AsText=$true}
    @{Cell='E4'; Value='-0.42%'; AsText=$true}
    @{Cell='G4'; Value='11'; AsText=$true}
    @{Cell='D5'; Value='0.06087'; AsText=$true}
    @{Cell='E5'; Value='2.38%'; AsText=$true}
    @{Cell='G5'; Value='11'; AsText=$true}
    @{Cell='D6'; Value='6.670'; AsText=$true}
    @{Cell='E6'; Value='0.75%'; AsText=$true}
    @{Cell='G6'; Value='11'; AsText=$true}
    @{Cell='B7'; Value='MXToken'; AsText=$false}
    @{Cell='C7'; Value='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; AsText=$false}
    @{Cell='D7'; Value='0.8469'; AsText=$true}
    @{Cell='E7'; Value='-0.36%'; AsText=$true}
    @{Cell='G7'; Value='11'; AsText=$true}
    @{Cell='B8'; Value='FTXToken'; AsText=$false}
    @{Cell='C8'; Value='https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; AsText=$false}
    @{Cell='D8'; Value='0.9297'; AsText=$true}
    @{Cell='E8'; Value='0.75%'; AsText=$true}
    @{Cell='G8'; Value='11'; AsText=$true}
    @{Cell='B9'; Value='WazirX'; AsText=$false}
    @{Cell='C9'; Value='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; AsText=$false}
    @{Cell='D9'; Value='0.1404'; AsText=$true}
    @{Cell='E9'; Value='1.73%'; AsText=$true}
    @{Cell='G9'; Value='11'; AsText=$true}
    @{Cell='B10'; Value='LiechtensteinCryptoassetsExchange'; AsText=$false}
    @{Cell='C10'; Value='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; AsText=$false}
    @{Cell='D10'; Value='0.04879'; AsText=$true}
    @{Cell='E10'; Value='15.74%'; AsText=$true}
    @{Cell='G10'; Value='11'; AsText=$true}
    @{Cell='B11'; Value='MandalaExchangeToken'; AsText=$false}
    @{Cell='C11'; Value='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; AsText=$false}
    @{Cell='D11'; Value='0.07099'; AsText=$true}
    @{Cell='E11'; Value='1.32%'; AsText=$true}
    @{Cell='G11'; Value='11'; AsText=$true}
    @{Cell='B12'; Value='BitrueCoin'; AsText=$false}
    @{Cell='C12'; Value='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; AsText=$false}
    @{Cell='D12'; Value='0.03080'; AsText=$true}
    @{Cell='E12'; Value='0.89%'; AsText=$true}
    @{Cell='G12'; Value='11'; AsText=$true}
    @{Cell='B13'; Value='BitMartToken'; AsText=$false}
    @{Cell='C13'; Value='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; AsText=$false}
    @{Cell='D13'; Value='0.09069'; AsText=$true}
    @{Cell='E13'; Value='-0.33%'; AsText=$true}
    @{Cell='G13'; Value='11'; AsText=$true}
    @{Cell='B14'; Value='BitForexToken'; AsText=$false}
    @{Cell='C14'; Value='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; AsText=$false}
    @{Cell='D14'; Value='0.001536'; AsText=$true}
    @{Cell='E14'; Value='-0.04%'; AsText=$true}
    @{Cell='G14'; Value='11'; AsText=$true}
    @{Cell='B15'; Value='One'; AsText=$false}
    @{Cell='C15'; Value='https://coinranking.com/coin/6Lga5NiXX3rT+one-one'; AsText=$false}
    @{Cell='D15'; Value='0.0006084'; AsText=$true}
    @{Cell='E15'; Value='0.88%'; AsText=$true}
    @{Cell='G15'; Value='11'; AsText=$true}
    @{Cell='B16'; Value='TigerCash'; AsText=$false}
    @{Cell='C16'; Value='https://coinranking.com/coin/6hIn06L2+tigercash-tch'; AsText=$false}
    @{Cell='D16'; Value='0.006016'; AsText=$true}
    @{Cell='E16'; Value='-0.98%'; AsText=$true}
    @{Cell='G16'; Value='11'; AsText=$true}
    @{Cell='B17'; Value='LEO'; AsText=$false}
    @{Cell='C17'; Value='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; AsText=$false}
    @{Cell='D17'; Value='3.450'; AsText=$true}
    @{Cell='E17'; Value='-0.48%'; AsText=$true}
    @{Cell='G17'; Value='11'; AsText=$true}
    @{Cell='B18'; Value='GateToken'; AsText=$false}
    @{Cell='C18'; Value='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; AsText=$false}
    @{Cell='D18'; Value='3.146'; AsText=$true}
    @{Cell='G18'; Value='11'; AsText=$true}
    @{Cell='E19'; Value='-0.76%'; AsText=$true}
    @{Cell='G19'; Value='11'; AsText=$true}
    @{Cell='E20'; Value='2.67%'; AsText=$true}
    @{Cell='G20'; Value='11'; AsText=$true}
    @{Cell='E21'; Value='0.24%'; AsText=$true}
    @{Cell='G21'; Value='11'; AsText=$true}
    @{Cell='D22'; Value='4.089'; AsText=$true}
    @{Cell='E22'; Value='3.82%'; AsText=$true}
    @{Cell='G22'; Value='11'; AsText=$true}
    @{Cell='D23'; Value='0.04240'; AsText=$true}
    @{Cell='E23'; Value='-0.53%'; AsText=$true}
    @{Cell='G23'; Value='11'; AsText=$true}
    @{Cell='E24'; Value='0.03%'; AsText=$true}
    @{Cell='G24'; Value='11'; AsText=$true}
    @{Cell='E25'; Value='4.96%'; AsText=$true}
    @{Cell='G25'; Value='11'; AsText=$true}
    @{Cell='E26'; Value='-0.08%'; AsText=$true}
    @{Cell='G26'; Value='11'; AsText=$true}
    @{Cell='D27'; Value='0.0001576'; AsText=$true}
    @{Cell='E27'; Value='3.41%'; AsText=$true}
    @{Cell='G27'; Value='11'; AsText=$true}
    @{Cell='G28'; Value='11'; AsText=$true}
    @{Cell='G29'; Value='11'; AsText=$true}
    @{Cell='G30'; Value='11'; AsText=$true}
    @{Cell='G31'; Value='11'; AsText=$true}
    @{Cell='G32'; Value='11'; AsText=$true}
    @{Cell='G33'; Value='11'; AsText=$true}
    @{Cell='G34'; Value='11'; AsText=$true}
    @{Cell='G35'; Value='11'; AsText=$true}
    @{Cell='G36'; Value='11'; AsText=$true}
    @{Cell='G37'; Value='11'; AsText=$true}
    @{Cell='G38'; Value='11'; AsText=$true}
    @{Cell='G39'; Value='11'; AsText=$true}
    @{Cell='D40'; Value='0.03867'; AsText=$true}
    @{Cell='E40'; Value='2.41%'; AsText=$true}
    @{Cell='G40'; Value='11'; AsText=$true}
    @{Cell='E41'; Value='1.54%'; AsText=$true}
    @{Cell='G41'; Value='11'; AsText=$true}
    @{Cell='E42'; Value='-35.07%'; AsText=$true}
    @{Cell='G42'; Value='11'; AsText=$true}
    @{Cell='D43'; Value='0.01635'; AsText=$true}
    @{Cell='E43'; Value='18.04%'; AsText=$true}
    @{Cell='G43'; Value='11'; AsText=$true}
    @{Cell='E44'; Value='-4.69%'; AsText=$true}
    @{Cell='G44'; Value='11'; AsText=$true}
    @{Cell='D45'; Value='0.00005148'; AsText=$true}
    @{Cell='E45'; Value='-3.29%'; AsText=$true}
    @{Cell='G45'; Value='11'; AsText=$true}
    @{Cell='E46'; Value='-0.04%'; AsText=$true}
    @{Cell='G46'; Value='11'; AsText=$true}
    @{Cell='G47'; Value='11'; AsText=$true}
    @{Cell='E48'; Value='23.73%'; AsText=$true}
    @{Cell='G48'; Value='11'; AsText=$true}
    @{Cell='E49'; Value='-0.04%'; AsText=$true}
    @{Cell='G49'; Value='11'; AsText=$true}
    @{Cell='E50'; Value='-0.04%'; AsText=$true}
    @{Cell='G50'; Value='11'; AsText=$true}
    @{Cell='G51'; Value='11'; AsText=$true}
)

foreach ($edit in $edits) {
    if ($edit.AsText) {
        # Leading apostrophe forces Excel to store the numeric-looking
        # value as text (matches t="inlineStr"/shared-string cells in
        # the source workbook) instead of coercing it to a number.
        $ws.Range($edit.Cell).Value = "'" + $edit.Value
    } else {
        $ws.Range($edit.Cell).Value = $edit.Value
    }
}

Write-Output ("Applied " + $edits.Count + " cell edits")
